$wb = $excel.ActiveWorkbook

# --- Sheet "inputdisp" ---
$ws1 = $wb.Worksheets.Item("inputdisp")
$ws1.Range("G2").Value = 20
$ws1.Range("H9").Select()

# --- Sheet "endofpipe" ---
$ws2 = $wb.Worksheets.Item("endofpipe")
$ws2.Range("G2").Value = 6
$ws2.Range("G7").Select()

# --- Sheet "inputprices" ---
$ws3 = $wb.Worksheets.Item("inputprices")
$ws3.Range("B2").Value = 3
$ws3.Range("B3").Value = 0.1
$ws3.Range("B4").Value = 3
$ws3.Range("B5").Value = 3
$ws3.Range("D13").Select()

$ws2.Activate()
